$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.979.21"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "2.357.04"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.69"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.65"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +4.87%  "

# Row 9
$ws.Range("E9").Value = "  +2.67%  "

# Row 10
$ws.Range("E10").Value = "  +2.83%  "

# Row 11
$ws.Range("E11").Value = "  -2.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("E12").Value = "  -1.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.92"
$ws.Range("E13").Value = "  +1.41%  "

# Row 14
$ws.Range("D14").Value = "2.781.13"
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("D15").Value = "57.934.94"
$ws.Range("E15").Value = "  -0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  +1.44%  "

# Row 17
$ws.Range("D17").Value = "2.350.36"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.94"
$ws.Range("E18").Value = "  +2.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.31"
$ws.Range("E19").Value = "  +2.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.94"
$ws.Range("E20").Value = "  -1.76%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("E21").Value = "  +3.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.42"
$ws.Range("E23").Value = "  +2.60%  "

# Row 24
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").Value = "  -2.92%  "

# Row 27
$ws.Range("E27").Value = "  -6.52%  "

# Row 28
$ws.Range("E28").Value = "  -0.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.46"
$ws.Range("E29").Value = "  +1.32%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0739"
$ws.Range("E30").Value = "  +1.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.35"
$ws.Range("E32").Value = "  -0.51%  "

# Row 33
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("E34").Value = "  -4.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.16"
$ws.Range("E36").Value = "  -1.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  -2.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").Value = "  -3.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.407"
$ws.Range("E39").Value = "  +7.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.45"
$ws.Range("E40").Value = "  -3.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.74"
$ws.Range("E42").Value = "  -0.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0947"
$ws.Range("E43").Value = "  +2.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0515"
$ws.Range("E44").Value = "  +2.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.563"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.65"
$ws.Range("E46").Value = "  -3.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0221"
$ws.Range("E47").Value = "  +1.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.388"
$ws.Range("E48").Value = "  +1.54%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.09"
$ws.Range("E49").Value = "  +0.26%  "

# Row 50
$ws.Range("E50").Value = "  +0.68%  "

# Row 51
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.947"
$ws.Range("E51").Value = "  +0.06%  "
